$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Shared teams/coaches" will not only play at different times but also
# will not play each other: the "Not Same Time As" note for 105-IRISH-WEBSTER 1
# ("Irish 2, Bears 3") is replaced by a Team ID reference (101) instead of text.
$ws.Range("I7").Value = 101

# Move the active selection cursor to I10 (cosmetic - cursor position only).
$ws.Range("I10").Select()
